# Auto-generated update of market-price-derived leve profit columns (H-N)
# Source data refreshed by the scheduled market-board runner; spreadsheet formulas
# are static snapshots (no formulas in these cells), so each touched cell is rewritten
# directly with its new value. A few cells go from populated -> blank (divide-by-zero
# guarded entries) and one goes from blank -> populated; those use ClearContents()/Value
# respectively so the serialized workbook matches exactly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H113").Value = 13312.3
$ws.Range("I113").Value = 7069.6
$ws.Range("J113").Value = 19555
$ws.Range("K113").Value = 7069.6
$ws.Range("L113").Value = 19555
$ws.Range("M113").Value = -3815.6
$ws.Range("N113").Value = -26063
$ws.Range("H116").Value = 3188.3333
$ws.Range("I116").Value = 2626.6667
$ws.Range("J116").Value = 3750
$ws.Range("K116").Value = 2626.6667
$ws.Range("L116").Value = 3750
$ws.Range("M116").Value = 815.3332999999998
$ws.Range("N116").Value = -10634
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()  # was -39820
$ws.Range("H131").Value = 7694.32
$ws.Range("I131").Value = 1202.6666
$ws.Range("K131").Value = 3607.9998
$ws.Range("M131").Value = 1432.0002
$ws.Range("H138").Value = 8067348
$ws.Range("I138").Value = 1319.579
$ws.Range("J138").Value = 11631407
$ws.Range("K138").Value = 3958.737
$ws.Range("L138").Value = 34894221
$ws.Range("M138").Value = 1181.263
$ws.Range("N138").Value = -34904501

$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 6758.018
$ws.Range("I32").Value = 3908.5107
$ws.Range("K32").Value = 3908.5107
$ws.Range("M32").Value = -3621.5107
$ws.Range("H74").Value = 8501.23
$ws.Range("I74").Value = 1819.1666
$ws.Range("J74").Value = 23535.875
$ws.Range("K74").Value = 1819.1666
$ws.Range("L74").Value = 23535.875
$ws.Range("M74").Value = -945.1666
$ws.Range("N74").Value = -25283.875
$ws.Range("H77").Value = 8501.23
$ws.Range("I77").Value = 1819.1666
$ws.Range("J77").Value = 23535.875
$ws.Range("K77").Value = 9095.833000000001
$ws.Range("L77").Value = 117679.375
$ws.Range("M77").Value = -4727.833000000001
$ws.Range("N77").Value = -126415.375
$ws.Range("H102").Value = 3238.5
$ws.Range("I102").Value = 2782.2
$ws.Range("K102").Value = 2782.2
$ws.Range("M102").Value = -1160.2
$ws.Range("H122").Value = 3149.8333
$ws.Range("I122").Value = 2999.5
$ws.Range("J122").Value = 3225
$ws.Range("K122").Value = 8998.5
$ws.Range("L122").Value = 9675
$ws.Range("M122").Value = -6548.5
$ws.Range("N122").Value = -14575

$ws = $wb.Worksheets("BSM")
$ws.Range("H20").Value = 2953.8044
$ws.Range("I20").Value = 2456.6206
$ws.Range("K20").Value = 2456.6206
$ws.Range("M20").Value = -2209.6206
$ws.Range("H107").Value = 1425.3572
$ws.Range("I107").Value = 1381.1538
$ws.Range("K107").Value = 1381.1538
$ws.Range("M107").Value = 538.8462

$ws = $wb.Worksheets("CRP")
$ws.Range("H31").Value = 45835.434
$ws.Range("I31").Value = 52170.8
$ws.Range("J31").Value = 3599.6667
$ws.Range("K31").Value = 52170.8
$ws.Range("L31").Value = 3599.6667
$ws.Range("M31").Value = -51875.8
$ws.Range("N31").Value = -4189.6667
$ws.Range("H34").Value = 45835.434
$ws.Range("I34").Value = 52170.8
$ws.Range("J34").Value = 3599.6667
$ws.Range("K34").Value = 52170.8
$ws.Range("L34").Value = 3599.6667
$ws.Range("M34").Value = -51968.8
$ws.Range("N34").Value = -4003.6667
$ws.Range("H105").Value = 1887.9412
$ws.Range("I105").Value = 636.875
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 636.875
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = 1110.125
$ws.Range("N105").Value = -6494
$ws.Range("H132").Value = 3867.7
$ws.Range("I132").Value = 3741.889
$ws.Range("K132").Value = 11225.667
$ws.Range("M132").Value = -8695.667000000001

$ws = $wb.Worksheets("CUL")
$ws.Range("H4").Value = 39789404
$ws.Range("J4").Value = 50019.5
$ws.Range("L4").Value = 150058.5
$ws.Range("N4").Value = -150282.5
$ws.Range("H68").Value = 8000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 8000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 24000
$ws.Range("M68").ClearContents()  # was -539
$ws.Range("N68").Value = -25622
$ws.Range("H71").Value = 8000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 8000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 72000
$ws.Range("M71").ClearContents()  # was 6
$ws.Range("N71").Value = -80112
$ws.Range("H122").Value = 1015.5714
$ws.Range("J122").Value = 1084.1666
$ws.Range("L122").Value = 9757.499400000001
$ws.Range("N122").Value = -14657.4994
$ws.Range("H129").Value = 1504.9166
$ws.Range("I129").Value = 1007.375
$ws.Range("K129").Value = 3022.125
$ws.Range("M129").Value = 1977.875
$ws.Range("H131").Value = 38881.52
$ws.Range("I131").Value = 501000
$ws.Range("K131").Value = 1503000
$ws.Range("M131").Value = -1497960

$ws = $wb.Worksheets("GSM")
$ws.Range("H38").Value = 22000
$ws.Range("J38").Value = 23000
$ws.Range("L38").Value = 23000
$ws.Range("N38").Value = -23926
$ws.Range("H52").Value = 41979.6
$ws.Range("J52").Value = 42474.5
$ws.Range("L52").Value = 42474.5
$ws.Range("N52").Value = -42992.5
$ws.Range("H102").Value = 25642460
$ws.Range("J102").Value = 1000000000
$ws.Range("L102").Value = 1000000000
$ws.Range("N102").Value = -1000003244
$ws.Range("H126").Value = 17477.066
$ws.Range("I126").Value = 24465.9
$ws.Range("J126").Value = 3499.4
$ws.Range("K126").Value = 73397.70000000001
$ws.Range("L126").Value = 10498.2
$ws.Range("M126").Value = -70927.70000000001
$ws.Range("N126").Value = -15438.2
$ws.Range("H134").Value = 72665.2
$ws.Range("J134").Value = 72665.2
$ws.Range("L134").Value = 217995.6
$ws.Range("N134").Value = -223065.6

$ws = $wb.Worksheets("LTW")
$ws.Range("H22").Value = 2008.7222
$ws.Range("I22").Value = 2395
$ws.Range("K22").Value = 2395
$ws.Range("M22").Value = -2100
$ws.Range("H27").Value = 2008.7222
$ws.Range("I27").Value = 2395
$ws.Range("K27").Value = 2395
$ws.Range("M27").Value = -2288
$ws.Range("H46").Value = 661.36365
$ws.Range("I46").Value = 595.8333
$ws.Range("J46").Value = 685.9375
$ws.Range("K46").Value = 595.8333
$ws.Range("L46").Value = 685.9375
$ws.Range("M46").Value = -407.8333
$ws.Range("N46").Value = -1061.9375
$ws.Range("H55").Value = 109.4
$ws.Range("I55").Value = 109.4
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 109.4
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = 63.59999999999999
$ws.Range("N55").ClearContents()  # was -606
$ws.Range("H122").Value = 4714.44
$ws.Range("I122").Value = 3881.1538
$ws.Range("J122").Value = 5617.1665
$ws.Range("K122").Value = 11643.4614
$ws.Range("L122").Value = 16851.4995
$ws.Range("M122").Value = -9193.4614
$ws.Range("N122").Value = -21751.4995
$ws.Range("H132").Value = 3562.6
$ws.Range("I132").Value = 2803.9
$ws.Range("J132").Value = 6597.4
$ws.Range("K132").Value = 8411.700000000001
$ws.Range("L132").Value = 19792.2
$ws.Range("M132").Value = -5881.700000000001
$ws.Range("N132").Value = -24852.2
$ws.Range("H136").Value = 3573.3447
$ws.Range("I136").Value = 3095.5
$ws.Range("K136").Value = 9286.5
$ws.Range("M136").Value = -6736.5

$ws = $wb.Worksheets("WVR")
$ws.Range("H132").Value = 1650.0714
$ws.Range("I132").Value = 1641.3334
$ws.Range("K132").Value = 4924.0002
$ws.Range("M132").Value = -2394.0002
$ws.Range("H133").Value = 75000
$ws.Range("J133").Value = 75000
$ws.Range("L133").Value = 75000
$ws.Range("N133").Value = -85120
$ws.Range("H137").Value = 49649.832
$ws.Range("J137").Value = 49649.832
$ws.Range("L137").Value = 49649.832
$ws.Range("N137").Value = -59849.832

